# "Update countries & provincias Spain"
#
# The sheet ("Pais") lists countries sorted descending by total cases
# (column B). This update refreshes the COVID-19 case counts for a
# handful of countries. Two of them (Ucrania / Oman) received big enough
# new totals to jump past their neighbours in the sort order, which
# cascades into the rows directly below them shifting down by one
# position (their data moves to the next row, like an inserted row).
# Everything else keeps its original row/position.
#
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: Austria - refreshed totals, no reordering ---
$ws.Cells.Item(18, 2).Value = 12675
$ws.Cells.Item(18, 3).Value = 36
$ws.Cells.Item(18, 5).Value = 8386

# --- Row 30: Chequia - refreshed totals, no reordering ---
$ws.Cells.Item(30, 2).Value = 5033
$ws.Cells.Item(30, 3).Value = 16
$ws.Cells.Item(30, 4).Value = 181
$ws.Cells.Item(30, 5).Value = 4761
$ws.Cells.Item(30, 6).Value = 103
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = 91

# --- Rows 54-57: Ucrania jumps ahead of Islandia/Singapur/Argelia ---
# Row 54 becomes Ucrania with its brand new totals; rows 55-57 each
# inherit the previous occupant's (old) data, shifted down by one.
$ws.Cells.Item(54, 1).Value = "Ucrania"
$ws.Cells.Item(54, 2).Value = 1668
$ws.Cells.Item(54, 3).Value = 206
$ws.Cells.Item(54, 4).Value = 35
$ws.Cells.Item(54, 5).Value = 1581
$ws.Cells.Item(54, 6).Value = 16
$ws.Cells.Item(54, 7).Value = 7
$ws.Cells.Item(54, 8).Value = 52

$ws.Cells.Item(55, 1).Value = "Islandia"
$ws.Cells.Item(55, 2).Value = 1586
$ws.Cells.Item(55, 4).Value = 559
$ws.Cells.Item(55, 5).Value = 1021
$ws.Cells.Item(55, 6).Value = 11

$ws.Cells.Item(56, 1).Value = "Singapur"
$ws.Cells.Item(56, 2).Value = 1481
$ws.Cells.Item(56, 4).Value = 377
$ws.Cells.Item(56, 5).Value = 1098
$ws.Cells.Item(56, 6).Value = 29
$ws.Cells.Item(56, 8).Value = 6

$ws.Cells.Item(57, 1).Value = "Argelia"
$ws.Cells.Item(57, 2).Value = 1468
$ws.Cells.Item(57, 4).Value = 113
$ws.Cells.Item(57, 5).Value = 1162
$ws.Cells.Item(57, 6).Value = 46
$ws.Cells.Item(57, 8).Value = 193
# Row 58 (Egipto) is unaffected and keeps its original values.

# --- Rows 90-94: Oman jumps ahead of Cuba/Burkina Faso/Albania/Taiwan ---
$ws.Cells.Item(90, 1).Value = "Oman"
$ws.Cells.Item(90, 2).Value = 419
$ws.Cells.Item(90, 3).Value = 48
$ws.Cells.Item(90, 4).Value = 72
$ws.Cells.Item(90, 5).Value = 345
$ws.Cells.Item(90, 6).Value = 3
$ws.Cells.Item(90, 8).Value = 2

$ws.Cells.Item(91, 1).Value = "Cuba"
$ws.Cells.Item(91, 2).Value = 396
$ws.Cells.Item(91, 4).Value = 27
$ws.Cells.Item(91, 5).Value = 358
$ws.Cells.Item(91, 6).Value = 15
$ws.Cells.Item(91, 8).Value = 11

$ws.Cells.Item(92, 1).Value = "Burkina Faso"
$ws.Cells.Item(92, 2).Value = 384
$ws.Cells.Item(92, 4).Value = 127
$ws.Cells.Item(92, 5).Value = 238
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 8).Value = 19

$ws.Cells.Item(93, 1).Value = "Albania"
$ws.Cells.Item(93, 2).Value = 383
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 131
$ws.Cells.Item(93, 5).Value = 230
$ws.Cells.Item(93, 6).Value = 7
$ws.Cells.Item(93, 8).Value = 22

$ws.Cells.Item(94, 1).Value = "Taiwan"
$ws.Cells.Item(94, 2).Value = 379
$ws.Cells.Item(94, 3).Value = 3
$ws.Cells.Item(94, 5).Value = 307
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 8).Value = 5
# Row 95 (Reunion) is unaffected and keeps its original values.

# --- Row 107: Vietnam - refreshed totals, no reordering ---
$ws.Cells.Item(107, 4).Value = 126
$ws.Cells.Item(107, 5).Value = 125

# --- Row 108: Montenegro - refreshed totals, no reordering ---
$ws.Cells.Item(108, 2).Value = 248
$ws.Cells.Item(108, 3).Value = 7
$ws.Cells.Item(108, 5).Value = 242

# --- Row 113: Islas Feroe - refreshed totals, no reordering ---
$ws.Cells.Item(113, 4).Value = 131
$ws.Cells.Item(113, 5).Value = 53
$ws.Cells.Item(113, 6).Value = 1

# --- Footer timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Abril de 2020 a las 08:52"
